$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF (Wins, Losses, Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of the existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2-51 with Wins=90, Losses=72, Ties=0
$ws.Range("AD2:AD51").Value = 90
$ws.Range("AE2:AE51").Value = 72
$ws.Range("AF2:AF51").Value = 0
